$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 90.36280833333332
$ws.Range("H2").Value = 271.088425
$ws.Range("I2").Value = 0.03168888268931816
$ws.Range("J2").Value = 0.03168888268931816
$ws.Range("M2").Value = 8.554479333333333
$ws.Range("N2").Value = 25.663438
$ws.Range("O2").Value = 0.1655051910559175
$ws.Range("P2").Value = 0.1655051910559175
$ws.Range("Q2").Value = 773.0067763894609
$ws.Range("R2").Value = 6957.060987505149
$ws.Range("S2").Value = 0.005244674583844159
$ws.Range("T2").Value = 0.00524467458384416
$ws.Range("G3").Value = 90.36280833333332
$ws.Range("H3").Value = 271.088425
$ws.Range("I3").Value = 0.03168888268931816
$ws.Range("J3").Value = 0.03168888268931816
$ws.Range("M3").Value = 20.28486166666667
$ws.Range("N3").Value = 60.854585
$ws.Range("O3").Value = 0.392455200938143
$ws.Range("P3").Value = 0.392455200938143
$ws.Range("Q3").Value = 1832.99706685318
$ws.Range("R3").Value = 16496.97360167862
$ws.Range("S3").Value = 0.0124364668233416
$ws.Range("T3").Value = 0.0124364668233416
$ws.Range("G4").Value = 90.36280833333332
$ws.Range("H4").Value = 271.088425
$ws.Range("I4").Value = 0.03168888268931816
$ws.Range("J4").Value = 0.03168888268931816
$ws.Range("M4").Value = 5.037112666666666
$ws.Range("N4").Value = 15.111338
$ws.Range("O4").Value = 0.09745400763531942
$ws.Range("P4").Value = 0.09745400763531943
$ws.Range("Q4").Value = 455.1676464514055
$ws.Range("R4").Value = 4096.50881806265
$ws.Range("S4").Value = 0.003088208615559553
$ws.Range("T4").Value = 0.003088208615559554
$ws.Range("G5").Value = 90.36280833333332
$ws.Range("H5").Value = 271.088425
$ws.Range("I5").Value = 0.03168888268931816
$ws.Range("J5").Value = 0.03168888268931816
$ws.Range("M5").Value = 17.810622
$ws.Range("N5").Value = 53.431866
$ws.Range("O5").Value = 0.34458560037062
$ws.Range("P5").Value = 0.34458560037062
$ws.Range("Q5").Value = 1609.41782208345
$ws.Range("R5").Value = 14484.76039875105
$ws.Range("S5").Value = 0.01091953266657284
$ws.Range("T5").Value = 0.01091953266657285
$ws.Range("I6").Value = 0.8807096817347263
$ws.Range("J6").Value = 0.8807096817347263
$ws.Range("M6").Value = 8.554479333333333
$ws.Range("N6").Value = 25.663438
$ws.Range("O6").Value = 0.1655051910559175
$ws.Range("P6").Value = 0.1655051910559175
$ws.Range("Q6").Value = 21483.70324972784
$ws.Range("R6").Value = 193353.3292475506
$ws.Range("S6").Value = 0.1457620241403022
$ws.Range("T6").Value = 0.1457620241403022
$ws.Range("I7").Value = 0.8807096817347263
$ws.Range("J7").Value = 0.8807096817347263
$ws.Range("M7").Value = 20.28486166666667
$ws.Range("N7").Value = 60.854585
$ws.Range("O7").Value = 0.392455200938143
$ws.Range("P7").Value = 0.392455200938143
$ws.Range("R7").Value = 458490.2696874851
$ws.Range("S7").Value = 0.3456390951133699
$ws.Range("T7").Value = 0.3456390951133699
$ws.Range("I8").Value = 0.8807096817347263
$ws.Range("J8").Value = 0.8807096817347263
$ws.Range("M8").Value = 5.037112666666666
$ws.Range("N8").Value = 15.111338
$ws.Range("O8").Value = 0.09745400763531942
$ws.Range("P8").Value = 0.09745400763531943
$ws.Range("Q8").Value = 12650.19524267699
$ws.Range("R8").Value = 113851.7571840929
$ws.Range("S8").Value = 0.08582868804827574
$ws.Range("T8").Value = 0.08582868804827576
$ws.Range("I9").Value = 0.8807096817347263
$ws.Range("J9").Value = 0.8807096817347263
$ws.Range("M9").Value = 17.810622
$ws.Range("N9").Value = 53.431866
$ws.Range("O9").Value = 0.34458560037062
$ws.Range("P9").Value = 0.34458560037062
$ws.Range("Q9").Value = 44729.56247028252
$ws.Range("R9").Value = 402566.0622325428
$ws.Range("S9").Value = 0.3034798744327783
$ws.Range("T9").Value = 0.3034798744327784
$ws.Range("G10").Value = 240.3144276666667
$ws.Range("H10").Value = 720.9432830000001
$ws.Range("I10").Value = 0.08427466838777388
$ws.Range("J10").Value = 0.08427466838777387
$ws.Range("M10").Value = 8.554479333333333
$ws.Range("N10").Value = 25.663438
$ws.Range("O10").Value = 0.1655051910559175
$ws.Range("P10").Value = 0.1655051910559175
$ws.Range("Q10").Value = 2055.764804976328
$ws.Range("R10").Value = 18501.88324478695
$ws.Range("S10").Value = 0.01394789509269261
$ws.Range("T10").Value = 0.01394789509269261
$ws.Range("G11").Value = 240.3144276666667
$ws.Range("H11").Value = 720.9432830000001
$ws.Range("I11").Value = 0.08427466838777388
$ws.Range("J11").Value = 0.08427466838777387
$ws.Range("M11").Value = 20.28486166666667
$ws.Range("N11").Value = 60.854585
$ws.Range("O11").Value = 0.392455200938143
$ws.Range("P11").Value = 0.392455200938143
$ws.Range("Q11").Value = 4874.744921722507
$ws.Range("R11").Value = 43872.70429550256
$ws.Range("S11").Value = 0.03307403191611916
$ws.Range("T11").Value = 0.03307403191611916
$ws.Range("G12").Value = 240.3144276666667
$ws.Range("H12").Value = 720.9432830000001
$ws.Range("I12").Value = 0.08427466838777388
$ws.Range("J12").Value = 0.08427466838777387
$ws.Range("M12").Value = 5.037112666666666
$ws.Range("N12").Value = 15.111338
$ws.Range("O12").Value = 0.09745400763531942
$ws.Range("P12").Value = 0.09745400763531943
$ws.Range("Q12").Value = 1210.490847582517
$ws.Range("R12").Value = 10894.41762824265
$ws.Range("S12").Value = 0.008212904176526127
$ws.Range("T12").Value = 0.008212904176526127
$ws.Range("G13").Value = 240.3144276666667
$ws.Range("H13").Value = 720.9432830000001
$ws.Range("I13").Value = 0.08427466838777388
$ws.Range("J13").Value = 0.08427466838777387
$ws.Range("M13").Value = 17.810622
$ws.Range("N13").Value = 53.431866
$ws.Range("O13").Value = 0.34458560037062
$ws.Range("P13").Value = 0.34458560037062
$ws.Range("Q13").Value = 4280.149432317342
$ws.Range("R13").Value = 38521.34489085608
$ws.Range("S13").Value = 0.02903983720243597
$ws.Range("T13").Value = 0.02903983720243597
$ws.Range("G14").Value = 9.486482333333333
$ws.Range("H14").Value = 28.459447
$ws.Range("I14").Value = 0.003326767188181744
$ws.Range("J14").Value = 0.003326767188181744
$ws.Range("M14").Value = 8.554479333333333
$ws.Range("N14").Value = 25.663438
$ws.Range("O14").Value = 0.1655051910559175
$ws.Range("P14").Value = 0.1655051910559175
$ws.Range("Q14").Value = 81.15191706653177
$ws.Range("R14").Value = 730.367253598786
$ws.Range("S14").Value = 0.0005505972390785771
$ws.Range("T14").Value = 0.0005505972390785772
$ws.Range("G15").Value = 9.486482333333333
$ws.Range("H15").Value = 28.459447
$ws.Range("I15").Value = 0.003326767188181744
$ws.Range("J15").Value = 0.003326767188181744
$ws.Range("M15").Value = 20.28486166666667
$ws.Range("N15").Value = 60.854585
$ws.Range("O15").Value = 0.392455200938143
$ws.Range("P15").Value = 0.392455200938143
$ws.Range("Q15").Value = 192.4319818349439
$ws.Range("R15").Value = 1731.887836514495
$ws.Range("S15").Value = 0.001305607085312287
$ws.Range("T15").Value = 0.001305607085312287
$ws.Range("G16").Value = 9.486482333333333
$ws.Range("H16").Value = 28.459447
$ws.Range("I16").Value = 0.003326767188181744
$ws.Range("J16").Value = 0.003326767188181744
$ws.Range("M16").Value = 5.037112666666666
$ws.Range("N16").Value = 15.111338
$ws.Range("O16").Value = 0.09745400763531942
$ws.Range("P16").Value = 0.09745400763531943
$ws.Range("Q16").Value = 47.78448032334288
$ws.Range("R16").Value = 430.060322910086
$ws.Range("S16").Value = 0.0003242067949579939
$ws.Range("T16").Value = 0.0003242067949579939
$ws.Range("G17").Value = 9.486482333333333
$ws.Range("H17").Value = 28.459447
$ws.Range("I17").Value = 0.003326767188181744
$ws.Range("J17").Value = 0.003326767188181744
$ws.Range("M17").Value = 17.810622
$ws.Range("N17").Value = 53.431866
$ws.Range("O17").Value = 0.34458560037062
$ws.Range("P17").Value = 0.34458560037062
$ws.Range("Q17").Value = 168.960150948678
$ws.Range("R17").Value = 1520.641358538102
$ws.Range("S17").Value = 0.001146356068832886
$ws.Range("T17").Value = 0.001146356068832886
